# forests-scraped.xlsx update - 2025-10-28 12:18
#
# The "New" sheet holds the most-recently-scraped listings; once they've
# been reviewed they get appended to the bottom of "Previously added" and
# the "New" sheet is refreshed with the next scrape batch.
#
# This script:
#   1. Moves the 9 existing rows on "New" (rows 2-10) onto the end of
#      "Previously added" (new rows 210-218), recreating their hyperlinks.
#   2. Clears "New" back down to just the header row.
#   3. Fills "New" with the 3 freshly scraped listings, with their own
#      hyperlinks.

$wb = $excel.ActiveWorkbook
$wsPrev = $wb.Worksheets.Item("Previously added")
$wsNew  = $wb.Worksheets.Item("New")

# ---------------------------------------------------------------------
# Data that is currently on "New" (rows 2-10), about to be archived into
# "Previously added". Columns: link, price, districtText, areaText,
# cadastreText, date (Excel serial).
# ---------------------------------------------------------------------
$movingRows = @(
    @("https://www.ss.com/msg/lv/real-estate/wood/balvi-and-reg/ziguru-pag/dfxge.html", "9 000 €", "Balvi un raj.", "1 ha.", "38980010052", 45957.40625),
    @("https://www.ss.com/msg/lv/real-estate/wood/cesis-and-reg/raunas-pag/jjxpi.html", "25 000 €", "Cēsis un raj.", "3 ha.", "42760080028", 45957.49444444444),
    @("https://www.ss.com/msg/lv/real-estate/wood/jekabpils-and-reg/garsenes-pag/bxojlf.html", "57 000 €", "Jēkabpils un raj.", "13 ha.", "56620020013", 45955.50069444445),
    @("https://www.ss.com/msg/lv/real-estate/wood/kraslava-and-reg/kalniesu-pag/jxgdk.html", "50 000 €", "Krāslava un raj.", "10 ha.", "60680040608", 45955.45),
    @("https://www.ss.com/msg/lv/real-estate/wood/liepaja-and-reg/durbe/gdkpm.html", "125 000 €", "Liepāja un raj.", "35 ha.", "64270060039", 45955.42152777778),
    @("https://www.ss.com/msg/lv/real-estate/wood/preili-and-reg/livani/jlngb.html", "22 500 €", "Preiļi un raj.", "4 ha.", "76860070184", 45957.44236111111),
    @("https://www.ss.com/msg/lv/real-estate/wood/rezekne-and-reg/kaunatas-pag/npjbg.html", "13 000 €", "Rēzekne un raj.", "4 ha.", "78620090056", 45957.57430555555),
    @("https://www.ss.com/msg/lv/real-estate/wood/rezekne-and-reg/feimanu-pag/mhedc.html", "50 000 €", "Rēzekne un raj.", "11.83 ha.", "78520030192", 45957.46805555555),
    @("https://www.ss.com/msg/lv/real-estate/wood/talsi-and-reg/kulciema-pag/cdlld.html", "51 000 €", "Talsi un raj.", "12 ha.", "88640040079", 45957.40833333333)
)

# ---------------------------------------------------------------------
# Freshly scraped rows for "New" (this batch's 3 listings).
# ---------------------------------------------------------------------
$freshRows = @(
    @("https://www.ss.com/msg/lv/real-estate/wood/daugavpils-and-reg/medumu-pag/conxg.html", "85 000 €", "Daugavpils un raj.", "23.60 ha.", "", 45958.38055555556),
    @("https://www.ss.com/msg/lv/real-estate/wood/kraslava-and-reg/udrisu-pag/jokck.html", "40 000 €", "Krāslava un raj.", "5 ha.", "60960050185", 45957.61597222222),
    @("https://www.ss.com/msg/lv/real-estate/wood/ludza-and-reg/pasienes-pag/dmgxk.html", "16 000 €", "Ludza un raj.", "22 ha.", "", 45958.51388888889)
)

# ---------------------------------------------------------------------
# Step 1: append the moving rows to the bottom of "Previously added".
# ---------------------------------------------------------------------
$destRow = $wsPrev.UsedRange.Rows.Count()

foreach ($row in $movingRows) {
    $destRow = $destRow + 1

    $wsPrev.Cells.Item($destRow, 1).Value = $row[0]
    $wsPrev.Cells.Item($destRow, 2).Value = $row[1]
    $wsPrev.Cells.Item($destRow, 3).Value = $row[2]
    $wsPrev.Cells.Item($destRow, 4).Value = $row[3]
    $wsPrev.Cells.Item($destRow, 5).Value = $row[4]
    $wsPrev.Cells.Item($destRow, 6).Value = $row[5]

    # Recreate the hyperlink on column A (this also (re)sets A's style to
    # the workbook's built-in Hyperlink style, so re-apply the sheet's own
    # look immediately afterwards).
    $wsPrev.Hyperlinks.Add($wsPrev.Cells.Item($destRow, 1), $row[0])

    # Match the formatting of the row directly above (style s="3" for A,
    # s="4" for B:E, s="2" for F) by copying it down.
    $srcRange = $wsPrev.Range("A" + ($destRow - 1) + ":F" + ($destRow - 1))
    $srcRange.Copy()
    $wsPrev.Range("A" + $destRow + ":F" + $destRow).PasteSpecial(-4122)
}

# ---------------------------------------------------------------------
# Step 2: reset "New" - drop every hyperlink (and the relationships that
# back them) then delete the now-archived rows, leaving just the header.
# ---------------------------------------------------------------------
$wsNew.Range("A1").Hyperlinks.Delete()
$wsNew.Range("A5:A10").EntireRow.Delete()

# ---------------------------------------------------------------------
# Step 3: write the 3 freshly scraped rows into "New" starting at row 2,
# reusing the existing row-2 formatting for each.
# ---------------------------------------------------------------------
$r = 2
foreach ($row in $freshRows) {
    $wsNew.Cells.Item($r, 1).Value = $row[0]
    $wsNew.Cells.Item($r, 2).Value = $row[1]
    $wsNew.Cells.Item($r, 3).Value = $row[2]
    $wsNew.Cells.Item($r, 4).Value = $row[3]
    $wsNew.Cells.Item($r, 5).Value = $row[4]
    $wsNew.Cells.Item($r, 6).Value = $row[5]

    $wsNew.Hyperlinks.Add($wsNew.Cells.Item($r, 1), $row[0])

    $r = $r + 1
}

# Re-apply the original row formatting (style s="3"/"4"/"2") across the
# rewritten data rows, now that the hyperlink writes are done.
$wsNew.Range("A2:F2").Copy()
$wsNew.Range("A2:F" + ($r - 1)).PasteSpecial(-4122)

Write-Host "Forests data updated."
